$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "PL20"
$ws.Range("B4").Value = "'12"
$ws.Range("C4").Value = 123123
$ws.Range("D4").Value = "asdasda"
